$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A-G)
# Row 2
$ws.Cells.Item(2, 1).Value = 87394
$ws.Cells.Item(2, 2).Value = "Pedro Novaes"
$ws.Cells.Item(2, 3).Value = "Engenharia"
$ws.Cells.Item(2, 4).Value = "Viagem de negócios"
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 45104
$ws.Cells.Item(2, 7).Value = 5634.67

# Row 3
$ws.Cells.Item(3, 1).Value = 17717
$ws.Cells.Item(3, 2).Value = "Beatriz Moreira"
$ws.Cells.Item(3, 3).Value = "Recursos Humanos"
$ws.Cells.Item(3, 4).Value = "Problemas pessoais"
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 45086
$ws.Cells.Item(3, 7).Value = 6204.58

# Row 4
$ws.Cells.Item(4, 1).Value = 17117
$ws.Cells.Item(4, 2).Value = "Rebeca Rocha"
$ws.Cells.Item(4, 3).Value = "Vendas"
$ws.Cells.Item(4, 4).Value = "Viagem de negócios"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 45097
$ws.Cells.Item(4, 7).Value = 2784.86

# Row 5
$ws.Cells.Item(5, 1).Value = 39200
$ws.Cells.Item(5, 2).Value = "Marina Gonçalves"
$ws.Cells.Item(5, 3).Value = "Operações"
$ws.Cells.Item(5, 4).Value = "Doença"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 45079
$ws.Cells.Item(5, 7).Value = 6814

# Row 6
$ws.Cells.Item(6, 1).Value = 37532
$ws.Cells.Item(6, 2).Value = "Noah Pires"
$ws.Cells.Item(6, 3).Value = "Vendas"
$ws.Cells.Item(6, 4).Value = "Doença"
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 45090
$ws.Cells.Item(6, 7).Value = 12085.82

# Row 7
$ws.Cells.Item(7, 1).Value = 18737
$ws.Cells.Item(7, 2).Value = "Diego Barbosa"
$ws.Cells.Item(7, 3).Value = "Financeiro"
$ws.Cells.Item(7, 4).Value = "Problemas pessoais"
$ws.Cells.Item(7, 5).Value = 6
$ws.Cells.Item(7, 6).Value = 45092
$ws.Cells.Item(7, 7).Value = 6468.52

# Row 8
$ws.Cells.Item(8, 1).Value = 85427
$ws.Cells.Item(8, 2).Value = "Gabriela Monteiro"
$ws.Cells.Item(8, 3).Value = "Jurídico"
$ws.Cells.Item(8, 4).Value = "Outros"
$ws.Cells.Item(8, 5).Value = 8
$ws.Cells.Item(8, 6).Value = 45085
$ws.Cells.Item(8, 7).Value = 9566.09

# Row 9
$ws.Cells.Item(9, 1).Value = 4347
$ws.Cells.Item(9, 2).Value = "Isabella Jesus"
$ws.Cells.Item(9, 3).Value = "P&D"
$ws.Cells.Item(9, 4).Value = "Doença"
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = 45100
$ws.Cells.Item(9, 7).Value = 4057.58

# Row 10
$ws.Cells.Item(10, 1).Value = 19511
$ws.Cells.Item(10, 2).Value = "Luigi Duarte"
$ws.Cells.Item(10, 3).Value = "Jurídico"
$ws.Cells.Item(10, 4).Value = "Viagem de negócios"
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = 45082
$ws.Cells.Item(10, 7).Value = 4619.27

# Row 11
$ws.Cells.Item(11, 1).Value = 61261
$ws.Cells.Item(11, 2).Value = "Raquel Freitas"
$ws.Cells.Item(11, 3).Value = "Marketing"
$ws.Cells.Item(11, 4).Value = "Problemas pessoais"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 45093
$ws.Cells.Item(11, 7).Value = 12415.01

$wb.Save()
